# Add new columns I (I0) and J (IF) to the sheet, mirroring the style of
# the existing header row and populating the data rows 2-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style used by the other header cells (e.g. H1) onto the new
# header cells so they match (bold font, border, centered/top alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-27) ---
$iValues = @(1,1,1,2,5,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$jValues = @(4,6,6,6,8,4,2,6,4,5,6,6,6,5,4,5,5,5,6,5,3,4,6,5,3,2)

for ($idx = 0; $idx -lt 26; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
